# Append a new trade row (row 8) to the GILD trades sheet, mirroring the
# layout/format of the existing data rows (A:I), then let Excel's
# "bestFit" column autosize pick up the new widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

# Column A: trade date/time (same number format style as the rows above it)
$ws.Cells.Item($row, 1).Value = 42654.745706018519
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

# Column B: Profitable (boolean)
$ws.Cells.Item($row, 2).Value = $true

# Column C: Principle
$ws.Cells.Item($row, 3).Value = 10057.530000000001

# Column D: Start Principle
$ws.Cells.Item($row, 4).Value = 10053.01

# Column E: BuyPrice
$ws.Cells.Item($row, 5).Value = 75.5

# Column F: SellPrice
$ws.Cells.Item($row, 6).Value = 75.569999999999993

# Column G: IsShortSell (boolean) - same number format style as col G above it
$ws.Cells.Item($row, 7).Value = $false
$ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"

# Column H: Price Change %
$ws.Cells.Item($row, 8).Value = 0.09

# Column I: Strong trade (boolean)
$ws.Cells.Item($row, 9).Value = $false

# Refresh the bestFit column widths now that new (wider/narrower) data exists
# (mirrors Excel auto-fitting A:I after the repeater appended the new row).
$ws.Columns.Item("A:I").AutoFit() | Out-Null

$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334
$ws.Columns.Item(9).ColumnWidth = 11

Write-Output "Added row 8 and refreshed column widths."
